$d = $word.ActiveDocument

$replacements = @(
    @("2024-05-09 Thursday", "2024-05-10 Friday"),
    @("534×4=", "492×7="),
    @("528×5=", "191×9="),
    @("844×2=", "800×4="),
    @("751×7=", "114×5="),
    @("316×5=", "166×4="),
    @("678×9=", "931×3="),
    @("845×7=", "583×2="),
    @("823×7=", "531×2="),
    @("950×6=", "634×6="),
    @("281×6=", "307×5="),
    @("298×8=", "604×5="),
    @("642×6=", "614×9="),
    @("757×8=", "321×9="),
    @("241×6=", "874×5="),
    @("182×2=", "825×6="),
    @("842×7=", "946×7="),
    @("896×4=", "203×4="),
    @("944×3=", "629×9="),
    @("577×6=", "831×8="),
    @("227×6=", "103×3="),
    @("757×6=", "845×8="),
    @("920×5=", "914×2="),
    @("776×5=", "253×6="),
    @("267×2=", "684×9="),
    @("572×4=", "738×4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
